# "Fixing Bulk Operation Template"
# Rename the sheet tab from "Create Item Group" to "Create Group".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Create Group"
